# Outstandings.xlsx - "Add files via upload" edit
#
# The Sale 22-23 sheet had three outdated ledger entries (old rows 8-10:
# invoice refs b23-24MQ114 / b23-24MQ205 / b23-24MQ207) removed, which also
# drops those three now-unused strings from the shared-strings table and
# renumbers every row below them. The running-total formula in F7 (which
# used to roll up through the deleted rows) is rewritten to sum just the
# remaining E5:E7 entries.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")
$ws2 = $wb.Worksheets.Item("Sale 22-23")

# Remove the three obsolete rows; everything beneath shifts up and Excel
# renumbers/relinks the dependent formulas automatically.
[void]$ws2.Rows("8:10").Delete()

# The subtotal that used to cover E5:E10 now only needs to cover E5:E7.
$ws2.Range("F7").Formula = "=E5+E6+E7"

# Restore the remembered cell selections for each sheet (captured at the
# time the file was last saved). Select Sale 22-23's cell first, then
# Purchase 22-23's, so Purchase 22-23 ends up the active/visible tab again.
[void]$ws2.Range("B31").Select()
[void]$ws1.Range("G33").Select()
